$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 149.5
$ws.Range("I4").Value = 149.5
$ws.Range("K4").Value = 149.5
$ws.Range("M4").Value = -35.5

$ws.Range("H5").Value = 173.1
$ws.Range("I5").Value = 188.5
$ws.Range("K5").Value = 188.5
$ws.Range("M5").Value = -73.5

$ws.Range("H9").Value = 80.333336
$ws.Range("J9").Value = 127.75
$ws.Range("L9").Value = 127.75
$ws.Range("N9").Value = -465.75

$ws.Range("H19").Value = 569.5
$ws.Range("J19").Value = 767
$ws.Range("L19").Value = 767
$ws.Range("N19").Value = -1117

$ws.Range("H39").Value = 394.86667
$ws.Range("I39").Value = 93.125
$ws.Range("J39").Value = 739.7143
$ws.Range("K39").Value = 279.375
$ws.Range("L39").Value = 2219.1429
$ws.Range("M39").Value = 16.625
$ws.Range("N39").Value = -2811.1429

$ws.Range("H53").Value = 382
$ws.Range("I53").Value = 503.66666
$ws.Range("J53").Value = 163
$ws.Range("K53").Value = 503.66666
$ws.Range("L53").Value = 163
$ws.Range("M53").Value = 133.33334
$ws.Range("N53").Value = -1437

$ws.Range("H76").Value = 4458.25
$ws.Range("I76").Value = 3333
$ws.Range("J76").Value = 4833.3335
$ws.Range("K76").Value = 3333
$ws.Range("L76").Value = 4833.3335
$ws.Range("M76").Value = -3018
$ws.Range("N76").Value = -5463.3335

$ws.Range("H79").Value = 4458.25
$ws.Range("I79").Value = 3333
$ws.Range("J79").Value = 4833.3335
$ws.Range("K79").Value = 3333
$ws.Range("L79").Value = 4833.3335
$ws.Range("M79").Value = -2241
$ws.Range("N79").Value = -7017.3335

$ws.Range("H113").Value = 2063.625
$ws.Range("I113").Value = 2063.625
$ws.Range("K113").Value = 2063.625
$ws.Range("M113").Value = 1190.375

$ws.Range("H116").Value = 7000
$ws.Range("J116").Value = 7000
$ws.Range("L116").Value = 7000
$ws.Range("N116").Value = -13884

$ws.Range("H137").Value = 3023.4167
$ws.Range("J137").Value = 3332.889
$ws.Range("L137").Value = 9998.667000000001
$ws.Range("N137").Value = -15098.667

$ws.Range("H138").Value = 2040.8334
$ws.Range("I138").Value = 1311.25
$ws.Range("J138").Value = 3500
$ws.Range("K138").Value = 3933.75
$ws.Range("L138").Value = 10500
$ws.Range("M138").Value = 1206.25
$ws.Range("N138").Value = -20780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4500
$ws.Range("I63").Value = 6000
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -5314
$ws.Range("N63").Value = -4372

$ws.Range("H66").Value = 4500
$ws.Range("I66").Value = 6000
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 30000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -26568
$ws.Range("N66").Value = -21864

$ws.Range("H74").Value = 6174.1665
$ws.Range("I74").Value = 5790
$ws.Range("K74").Value = 5790
$ws.Range("M74").Value = -4916

$ws.Range("H77").Value = 6174.1665
$ws.Range("I77").Value = 5790
$ws.Range("K77").Value = 28950
$ws.Range("M77").Value = -24582

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H101").Value = 34000
$ws.Range("J101").Value = 34000
$ws.Range("L101").Value = 34000
$ws.Range("N101").Value = -40490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5469.659
$ws.Range("I31").Value = 2992.8948
$ws.Range("K31").Value = 2992.8948
$ws.Range("M31").Value = -2697.8948

$ws.Range("H34").Value = 5469.659
$ws.Range("I34").Value = 2992.8948
$ws.Range("K34").Value = 2992.8948
$ws.Range("M34").Value = -2790.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 96.57143000000001
$ws.Range("I33").Value = 116.333336
$ws.Range("J33").Value = 81.75
$ws.Range("K33").Value = 698.000016
$ws.Range("L33").Value = 490.5
$ws.Range("M33").Value = -415.000016
$ws.Range("N33").Value = -1056.5

$ws.Range("H34").Value = 1851.258
$ws.Range("I34").Value = 166.2
$ws.Range("J34").Value = 2175.3076
$ws.Range("K34").Value = 498.6
$ws.Range("L34").Value = 6525.9228
$ws.Range("M34").Value = -414.6
$ws.Range("N34").Value = -6693.9228

$ws.Range("H68").Value = 583.1667
$ws.Range("J68").Value = 599.8
$ws.Range("L68").Value = 1799.4
$ws.Range("N68").Value = -3421.4

$ws.Range("H71").Value = 583.1667
$ws.Range("J71").Value = 599.8
$ws.Range("L71").Value = 5398.2
$ws.Range("N71").Value = -13510.2

$ws.Range("H113").Value = 1671.1333
$ws.Range("J113").Value = 1984.4
$ws.Range("L113").Value = 5953.200000000001
$ws.Range("N113").Value = -10293.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4266.1665
$ws.Range("I80").Value = 4124.25
$ws.Range("K80").Value = 4124.25
$ws.Range("M80").Value = -3126.25

$ws.Range("H83").Value = 4266.1665
$ws.Range("I83").Value = 4124.25
$ws.Range("K83").Value = 20621.25
$ws.Range("M83").Value = -15629.25

$ws.Range("H98").Value = 9170
$ws.Range("J98").Value = 9170
$ws.Range("L98").Value = 9170
$ws.Range("N98").Value = -15160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1330
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""

$ws.Range("H22").Value = 2502.75
$ws.Range("J22").Value = 3142.8572
$ws.Range("L22").Value = 3142.8572
$ws.Range("N22").Value = -3732.8572

$ws.Range("H27").Value = 2502.75
$ws.Range("J27").Value = 3142.8572
$ws.Range("L27").Value = 3142.8572
$ws.Range("N27").Value = -3356.8572

$ws.Range("H116").Value = 252000
$ws.Range("J116").Value = 252000
$ws.Range("L116").Value = 252000
$ws.Range("N116").Value = -261178

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1424.5
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 1432.6666
$ws.Range("K113").Value = 4200
$ws.Range("L113").Value = 4297.9998
$ws.Range("M113").Value = -2030
$ws.Range("N113").Value = -8637.9998

$ws.Range("H126").Value = 5937.75
$ws.Range("I126").Value = 3100.6
$ws.Range("J126").Value = 7964.2856
$ws.Range("K126").Value = 9301.799999999999
$ws.Range("L126").Value = 23892.8568
$ws.Range("M126").Value = -6831.799999999999
$ws.Range("N126").Value = -28832.8568

$ws.Range("H136").Value = 2508.1462
$ws.Range("I136").Value = 1708.1613
$ws.Range("K136").Value = 5124.4839
$ws.Range("M136").Value = -2574.4839
